$d = $word.ActiveDocument

$pairs = @(
    @("18×61=1098", "96×47=4512"),
    @("37×53=1961", "13×58=754"),
    @("34×18=612", "81×82=6642"),
    @("46×12=552", "95×48=4560"),
    @("91×58=5278", "80×49=3920"),
    @("80×87=6960", "30×65=1950"),
    @("81×50=4050", "86×66=5676"),
    @("61×44=2684", "33×38=1254"),
    @("52×41=2132", "14×41=574"),
    @("49×58=2842", "16×97=1552"),
    @("28×72=2016", "35×55=1925"),
    @("73×49=3577", "95×83=7885"),
    @("60×40=2400", "71×33=2343"),
    @("43×65=2795", "25×72=1800"),
    @("36×28=1008", "74×77=5698"),
    @("53×19=1007", "70×68=4760"),
    @("29×80=2320", "32×94=3008"),
    @("67×54=3618", "54×23=1242"),
    @("35×71=2485", "43×11=473"),
    @("38×78=2964", "80×20=1600"),
    @("53×80=4240", "25×45=1125"),
    @("93×27=2511", "12×57=684"),
    @("16×91=1456", "72×92=6624"),
    @("88×25=2200", "72×50=3600"),
    @("96×95=9120", "96×87=8352")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
